$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly row: copy the current row 5 (week of 2021-10-08) down to row 6
# before row 5 is overwritten with the new week's data.
$ws.Cells.Item(6, 1).Value = $ws.Cells.Item(5, 1).Value()
$ws.Cells.Item(6, 2).Value = $ws.Cells.Item(5, 2).Value()
$ws.Cells.Item(6, 3).Value = $ws.Cells.Item(5, 3).Value()
$ws.Cells.Item(6, 4).NumberFormat = $ws.Cells.Item(5, 4).NumberFormat()
$ws.Cells.Item(6, 4).Value = $ws.Cells.Item(5, 4).Value()
$ws.Cells.Item(6, 5).Value = $ws.Cells.Item(5, 5).Value()
$ws.Cells.Item(6, 6).Value = $ws.Cells.Item(5, 6).Value()
$ws.Cells.Item(6, 7).Value = $ws.Cells.Item(5, 7).Value()
$ws.Cells.Item(6, 8).Value = $ws.Cells.Item(5, 8).Value()
$ws.Cells.Item(6, 9).Value = $ws.Cells.Item(5, 9).Value()
$ws.Cells.Item(6, 10).Value = $ws.Cells.Item(5, 10).Value()
$ws.Cells.Item(6, 11).Value = $ws.Cells.Item(5, 11).Value()
$ws.Cells.Item(6, 12).Value = $ws.Cells.Item(5, 12).Value()
$ws.Cells.Item(6, 13).Value = $ws.Cells.Item(5, 13).Value()
$ws.Cells.Item(6, 14).Value = $ws.Cells.Item(5, 14).Value()
$ws.Cells.Item(6, 15).Value = $ws.Cells.Item(5, 15).Value()
$ws.Cells.Item(6, 16).Value = $ws.Cells.Item(5, 16).Value()
$ws.Cells.Item(6, 17).Value = $ws.Cells.Item(5, 17).Value()
$ws.Cells.Item(6, 18).Value = $ws.Cells.Item(5, 18).Value()

# Row 5 now holds the newest week's report (2021-11-10).
$ws.Cells.Item(5, 4).Value = 44510
$ws.Cells.Item(5, 10).Value = 600
$ws.Cells.Item(5, 11).Value = 1300
$ws.Cells.Item(5, 12).Value = 1400
$ws.Cells.Item(5, 13).Value = 1350
$ws.Cells.Item(5, 16).Value = 1350
